$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix another movement bug: update a handful of map cells between
# "open" and "wall" (shared string values 0/1 respectively).

$ws.Range("C8").Value = "open"

$ws.Range("B9").Value  = "wall"
$ws.Range("D9").Value  = "wall"

$ws.Range("B10").Value = "wall"
$ws.Range("D10").Value = "wall"

$ws.Range("B11").Value = "wall"
$ws.Range("D11").Value = "wall"

$ws.Range("B12").Value = "wall"
$ws.Range("D12").Value = "wall"

$ws.Range("B13").Value = "wall"
$ws.Range("D13").Value = "wall"
$ws.Range("E13").Value = "wall"
$ws.Range("F13").Value = "wall"

$ws.Range("B14").Value = "wall"

$ws.Range("B15").Value = "wall"
$ws.Range("C15").Value = "wall"
$ws.Range("D15").Value = "wall"
$ws.Range("E15").Value = "wall"
$ws.Range("F15").Value = "wall"
